# "get data to model" - the last sample row (row 60) on the meter-reading
# sheet gets its data wiped out (A:C lose the cell entirely, D:E keep their
# style but lose the value), which also drops the now-orphaned shared
# string "2015/10/27 10:43:28" from the workbook's string table once it is
# no longer referenced anywhere.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Wipe out the last row's data (A60:E60). D60/E60 keep their number-format
# style (s="1") but end up with no value, same as A60:C60 which end up
# with no cell content at all.
$ws.Range("A60:E60").ClearContents() | Out-Null

# Scroll the window down so row 60 is visible near the bottom, and leave
# the new empty row selected - matches the author's view state when they
# made the edit.
$excel.ActiveWindow.ScrollRow = 55
$ws.Range("A60:G60").Select() | Out-Null
